# Update the quarterly report: drop the oldest quarter column (1399/06)
# and shift all subsequent quarters one column to the left, appending the
# newest quarter (1401/12) in the last column, with its figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header rows (quarter labels) -------------------------------------
$headers = @(
    "فصل سوم منتهی به 1399/09",
    "فصل چهارم منتهی به 1399/12",
    "فصل اول منتهی به 1400/03",
    "فصل دوم منتهی به 1400/06",
    "فصل سوم منتهی به 1400/09",
    "فصل چهارم منتهی به 1400/12",
    "فصل اول منتهی به 1401/03",
    "فصل دوم منتهی به 1401/06",
    "فصل سوم منتهی به 1401/09",
    "فصل چهارم منتهی به 1401/12"
)
$cols = @("E","F","G","H","I","J","K","L","M","N")

for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "8").Value = $headers[$i]
    $ws.Range($cols[$i] + "24").Value = $headers[$i]
}

# --- Data rows (shift one column left, append new quarter value) ------
$row17 = @(202577, 430524, 480105, 571769, 437670, 306288, 682834, 1266066, 1012820, 657108)
$row19 = @(170724, 7878465, 560959, 522804, 1704659, 5683364, 1348488, 1870510, -406430, 7582329)
$row20 = @(373301, 8308989, 1041064, 1094573, 2142329, 5989652, 2031322, 3136576, 606390, 8239437)
$row26 = @(789, 790, 790, 775, 775, 756, 756, 748, 748, 739)
$row27 = @(1613, 1556, 1556, 1535, 1535, 1520, 1520, 1504, 1504, 1487)

for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "17").Value = $row17[$i]
    $ws.Range($cols[$i] + "19").Value = $row19[$i]
    $ws.Range($cols[$i] + "20").Value = $row20[$i]
    $ws.Range($cols[$i] + "26").Value = $row26[$i]
    $ws.Range($cols[$i] + "27").Value = $row27[$i]
}

$wb.Save()
